# Auto-generated Excel COM-interop edit script
# Applies numeric updates to ALC, ARM, BSM, CRP, CUL, LTW, WVR sheets
# per the target diff (Pandaemonium_Profits workbook scheduled-runner update).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 57767.277
$ws.Range("I86").Value = 168200.5
$ws.Range("J86").Value = 2550.6667
$ws.Range("K86").Value = 168200.5
$ws.Range("L86").Value = 2550.6667
$ws.Range("M86").Value = -167077.5
$ws.Range("N86").Value = -4796.6667

$ws.Range("H88").Value = 3000
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 3000
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 3000
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -3812

$ws.Range("H89").Value = 57767.277
$ws.Range("I89").Value = 168200.5
$ws.Range("J89").Value = 2550.6667
$ws.Range("K89").Value = 841002.5
$ws.Range("L89").Value = 12753.3335
$ws.Range("M89").Value = -835386.5
$ws.Range("N89").Value = -23985.3335

$ws.Range("H91").Value = 3000
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 3000
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 3000
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -5808

$ws.Range("H135").Value = 41667410
$ws.Range("I135").Value = 16129804
$ws.Range("J135").Value = 200000560
$ws.Range("K135").Value = 145168236
$ws.Range("L135").Value = 1800005040
$ws.Range("M135").Value = -145165701
$ws.Range("N135").Value = -1800010110

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 9213.9
$ws.Range("I61").Value = 6161.324
$ws.Range("J61").Value = 17902
$ws.Range("K61").Value = 6161.324
$ws.Range("L61").Value = 17902
$ws.Range("M61").Value = -5949.324
$ws.Range("N61").Value = -18326

$ws.Range("H86").Value = 40314
$ws.Range("J86").Value = 40314
$ws.Range("L86").Value = 40314
$ws.Range("N86").Value = -42686

$ws.Range("H88").Value = 20733.334
$ws.Range("I88").Value = 20733.334
$ws.Range("K88").Value = 20733.334
$ws.Range("M88").Value = -20327.334

$ws.Range("H89").Value = 40314
$ws.Range("J89").Value = 40314
$ws.Range("L89").Value = 120942
$ws.Range("N89").Value = -132798

$ws.Range("H91").Value = 20733.334
$ws.Range("I91").Value = 20733.334
$ws.Range("K91").Value = 20733.334
$ws.Range("M91").Value = -19329.334

$ws.Range("H132").Value = 3868.322
$ws.Range("I132").Value = 1700.7941
$ws.Range("J132").Value = 6816.16
$ws.Range("K132").Value = 5102.3823
$ws.Range("L132").Value = 20448.48
$ws.Range("M132").Value = -2572.3823
$ws.Range("N132").Value = -25508.48

$ws.Range("H136").Value = 9213.9
$ws.Range("I136").Value = 6161.324
$ws.Range("J136").Value = 17902
$ws.Range("K136").Value = 18483.972
$ws.Range("L136").Value = 53706
$ws.Range("M136").Value = -15933.972
$ws.Range("N136").Value = -58806

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1892.1052
$ws.Range("I86").Value = 1902.7778
$ws.Range("J86").Value = 1700
$ws.Range("K86").Value = 1902.7778
$ws.Range("L86").Value = 1700
$ws.Range("M86").Value = -779.7778000000001
$ws.Range("N86").Value = -3946

$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()

$ws.Range("H89").Value = 1892.1052
$ws.Range("I89").Value = 1902.7778
$ws.Range("J89").Value = 1700
$ws.Range("K89").Value = 9513.889000000001
$ws.Range("L89").Value = 8500
$ws.Range("M89").Value = -3897.889000000001
$ws.Range("N89").Value = -19732

$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6552.0347
$ws.Range("I31").Value = 9226.786
$ws.Range("J31").Value = 4055.6
$ws.Range("K31").Value = 9226.786
$ws.Range("L31").Value = 4055.6
$ws.Range("M31").Value = -8931.786
$ws.Range("N31").Value = -4645.6

$ws.Range("H34").Value = 6552.0347
$ws.Range("I34").Value = 9226.786
$ws.Range("J34").Value = 4055.6
$ws.Range("K34").Value = 9226.786
$ws.Range("L34").Value = 4055.6
$ws.Range("M34").Value = -9024.786
$ws.Range("N34").Value = -4459.6

$ws.Range("H45").Value = 12000
$ws.Range("I45").Value = 9000
$ws.Range("J45").Value = 15000
$ws.Range("K45").Value = 9000
$ws.Range("L45").Value = 15000
$ws.Range("M45").Value = -8407
$ws.Range("N45").Value = -16186

$ws.Range("H97").Value = 70540
$ws.Range("I97").Value = 70080
$ws.Range("J97").Value = 71000
$ws.Range("K97").Value = 70080
$ws.Range("L97").Value = 71000
$ws.Range("M97").Value = -69089
$ws.Range("N97").Value = -72982

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 713.5714
$ws.Range("I113").Value = 710.4167
$ws.Range("J113").Value = 732.5
$ws.Range("K113").Value = 2131.2501
$ws.Range("L113").Value = 2197.5
$ws.Range("M113").Value = 38.7498999999998
$ws.Range("N113").Value = -6537.5

$ws.Range("H122").Value = 998.9167
$ws.Range("J122").Value = 1042.5714
$ws.Range("L122").Value = 9383.142600000001
$ws.Range("N122").Value = -14283.1426

$ws.Range("H132").Value = 1674.1428
$ws.Range("J132").Value = 1434
$ws.Range("L132").Value = 12906
$ws.Range("N132").Value = -17966

$ws.Range("H136").Value = 2942.8262
$ws.Range("I136").Value = 1263.3334
$ws.Range("J136").Value = 3535.5881
$ws.Range("K136").Value = 3790.0002
$ws.Range("L136").Value = 10606.7643
$ws.Range("M136").Value = 1309.9998
$ws.Range("N136").Value = -20806.7643

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

$ws.Range("H88").Value = 19997.8
$ws.Range("I88").Value = 14950
$ws.Range("J88").Value = 40189
$ws.Range("K88").Value = 14950
$ws.Range("L88").Value = 40189
$ws.Range("M88").Value = -14522
$ws.Range("N88").Value = -41045

$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

$ws.Range("H91").Value = 19997.8
$ws.Range("I91").Value = 14950
$ws.Range("J91").Value = 40189
$ws.Range("K91").Value = 14950
$ws.Range("L91").Value = 40189
$ws.Range("M91").Value = -13468
$ws.Range("N91").Value = -43153

$ws.Range("H132").Value = 5068.7666
$ws.Range("I132").Value = 5639.7383
$ws.Range("J132").Value = 3736.5
$ws.Range("K132").Value = 16919.2149
$ws.Range("L132").Value = 11209.5
$ws.Range("M132").Value = -14389.2149
$ws.Range("N132").Value = -16269.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 941.94446
$ws.Range("I100").Value = 206.5
$ws.Range("K100").Value = 413
$ws.Range("M100").Value = 128
